# Updated cryptos list on Sat Sep 28 23:59:50 UTC 2024 with GitHub Actions
#
# This script applies the latest crypto price/volume snapshot to Sheet1.
# Columns: A=rank(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)
#
# All text-like values (prices, percentages, coin names, links) must stay
# plain text cells, matching the original inlineStr cell type. Excel's COM
# layer auto-coerces strings that look numeric (e.g. "1.00") into real
# numbers, so each target cell is forced to Text format before the value
# is assigned, then the cell style is reset back to "Normal" so no stray
# number-format style gets attached (keeping parity with the un-styled
# cells in the original workbook).
#
# NOTE: this runtime's PowerShell-like parser only binds positional
# function arguments (named "-Param value" args are silently ignored), so
# all helper calls below use positional args with $null placeholders for
# columns that do not change on a given row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

function Set-Row($Row, $B, $C, $D, $E) {
    if ($B) { Set-TextValue $ws.Cells.Item($Row, 2) $B }
    if ($C) { Set-TextValue $ws.Cells.Item($Row, 3) $C }
    if ($D) { Set-TextValue $ws.Cells.Item($Row, 4) $D }
    if ($E) { Set-TextValue $ws.Cells.Item($Row, 5) $E }
}

# --- Simple price / volume refreshes (no identity changes) ---------------

Set-Row 2  $null $null "65.917.83" "  +0.17%  "   # Bitcoin
Set-Row 3  $null $null "2.679.17"  "  -0.62%  "   # Ethereum
Set-Row 4  $null $null $null       "  -0.01%  "   # TetherUSD
Set-Row 5  $null $null "601.45"    "  -0.98%  "   # BNB
Set-Row 6  $null $null "156.91"    "  -0.52%  "   # Solana
Set-Row 7  $null $null "1.00"      "  -0.03%  "   # USDC
Set-Row 8  $null $null "0.615"     "  +4.39%  "   # XRP
Set-Row 9  $null $null "0.129"     "  +3.94%  "   # Dogecoin
Set-Row 10 $null $null "0.401"     "  -0.31%  "   # Cardano
Set-Row 11 $null $null "5.89"      "  -1.53%  "   # Toncoin
Set-Row 12 $null $null $null       "  -0.04%  "   # TRON
Set-Row 13 $null $null "29.34"     "  -3.66%  "   # Avalanche
Set-Row 14 $null $null "0.0000197" "  -3.12%  "   # ShibaInu
Set-Row 15 $null $null "3.158.35"  "  -0.78%  "   # WrappedliquidstakedEther2.0
Set-Row 16 $null $null "65.728.44" "  +0.14%  "   # WrappedBTC
Set-Row 17 $null $null "2.681.97"  "  -0.49%  "   # WrappedEther
Set-Row 18 $null $null "12.95"     "  +2.19%  "   # Chainlink
Set-Row 19 $null $null $null       "  -1.72%  "   # Polkadot
Set-Row 20 $null $null "7.61"      "  -0.11%  "   # Uniswap
Set-Row 21 $null $null "353.15"    "  -1.64%  "   # BitcoinCash
Set-Row 22 $null $null $null       "  +0.07%  "   # Dai
Set-Row 23 $null $null "69.96"     "  -1.60%  "   # Litecoin
Set-Row 24 $null $null "0.0000114" "  +6.57%  "   # PEPE
Set-Row 25 $null $null "9.70"      "  -1.86%  "   # InternetComputer(DFINITY)
Set-Row 26 $null $null $null       "  +1.66%  "   # SuiNetwork
Set-Row 27 $null $null $null       "  -2.65%  "   # Kaspa
Set-Row 28 $null $null $null       "  -4.88%  "   # Fetch.AI
Set-Row 29 $null $null "8.07"      "  -4.77%  "   # Aptos
Set-Row 30 $null $null "0.997"     "  -0.46%  "   # Binance-PegBSC-USD
Set-Row 31 $null $null "536.32"    "  -2.44%  "   # Bittensor
Set-Row 32 $null $null $null       "  -2.65%  "   # PancakeSwap
Set-Row 33 $null $null "1.78"      "  -2.26%  "   # ImmutableX

# --- Rows 34/35 swap identity: NEARProtocol <-> RenderToken ---------------

Set-Row 34 "RenderToken" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render" "6.54" "  -2.41%  "
Set-Row 35 "NEARProtocol" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" "5.56" "  +1.96%  "

# --- Simple refreshes continue -------------------------------------------

Set-Row 36 $null $null "0.424" "  -2.21%  "   # PolygonEcosystemToken
Set-Row 37 $null $null "20.59" "  -1.31%  "   # EthereumClassic
Set-Row 38 $null $null "1.00"  "  +0.02%  "   # FirstDigitalUSD

# --- Rows 39/40 swap identity: Monero <-> Stacks --------------------------

Set-Row 39 "Stacks" "https://coinranking.com/coin/mMPrMcB7+stacks-stx" "1.96" "  -1.43%  "
Set-Row 40 "Monero" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" "158.48" "  -3.33%  "

# --- Simple refreshes continue -------------------------------------------

Set-Row 41 $null $null $null    "  +0.03%  "   # USDe
Set-Row 42 $null $null "163.98" "  -3.89%  "   # Aave
Set-Row 43 $null $null "4.13"   "  -1.62%  "   # Filecoin
Set-Row 44 $null $null "2.37"   "  +4.65%  "   # dogwifhat
Set-Row 45 $null $null "0.0613" "  -1.12%  "   # Hedera
Set-Row 46 $null $null "22.88"  "  -3.18%  "   # InjectiveProtocol

# --- Rows 47/49 swap identity: BabyDogeCoin <-> Mantle (row 48 untouched) -

Set-Row 47 "Mantle" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" "0.641" "  -2.17%  "

Set-Row 48 $null $null "0.0259" "  -2.92%  "   # VeChain

Set-Row 49 "BabyDogeCoin" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" "0.0₆0265" "  +14.66%  "

# --- Final simple refreshes ------------------------------------------------

Set-Row 50 $null $null "20.20"  "  -3.69%  "   # EnergySwap
Set-Row 51 $null $null "0.0995" "  +0.10%  "   # Stellar
